# Adds the four new FAQ tag groups (law-why, law-calculate, law-time,
# law-place, law-salary) as new rows 58-94 on Sheet1, each row holding a
# "tag" value in column A and a matching training phrase in column B -
# mirroring the existing tag/question table already present in the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 58-64 (the "law-why" group) continue to carry the bold-less "tag"
# style (s=2) used by column A throughout the rest of the table, so copy
# that formatting down from the last existing row (A57:B57) first.
$fmtSource = $ws.Range("A57:B57")

$fmtSource.Copy()
$ws.Range("A58:B58").PasteSpecial(-4122)
$ws.Range("A58").Value = "law-why"
$ws.Range("B58").Value = "เหตุใดทำไมหน่วยงานรัฐถึงต้องจัดเก็บภาษี?"

$fmtSource.Copy()
$ws.Range("A59:B59").PasteSpecial(-4122)
$ws.Range("A59").Value = "law-why"
$ws.Range("B59").Value = "ทำไมหน่วยงานรัฐต้องเก็บภาษี"

$fmtSource.Copy()
$ws.Range("A60:B60").PasteSpecial(-4122)
$ws.Range("A60").Value = "law-why"
$ws.Range("B60").Value = "ทำไมต้องเก็บภาษี"

$fmtSource.Copy()
$ws.Range("A61:B61").PasteSpecial(-4122)
$ws.Range("A61").Value = "law-why"
$ws.Range("B61").Value = "หน่วยงานรัฐเก็บภาษีทำไม"

$fmtSource.Copy()
$ws.Range("A62:B62").PasteSpecial(-4122)
$ws.Range("A62").Value = "law-why"
$ws.Range("B62").Value = "หน่วยงานรัฐเก็บภาษีเพื่อ"

$fmtSource.Copy()
$ws.Range("A63:B63").PasteSpecial(-4122)
$ws.Range("A63").Value = "law-why"
$ws.Range("B63").Value = "หน่วยงานรัฐเก็บภาษีเพื่ออะไร"

$fmtSource.Copy()
$ws.Range("A64:B64").PasteSpecial(-4122)
$ws.Range("A64").Value = "law-why"
$ws.Range("B64").Value = "ทำไมต้องเก็บภาษี"

$ws.Range("A65").Value = "law-calculate"
$ws.Range("B65").Value = "ภาษีบุคคลธรรมดาคำนวณจากอะไรและคำนวนอย่างไร"

$ws.Range("A66").Value = "law-calculate"
$ws.Range("B66").Value = "คำนวณภาษีอย่างไร"

$ws.Range("A67").Value = "law-calculate"
$ws.Range("B67").Value = "คำนวณภาษี"

$ws.Range("A68").Value = "law-calculate"
$ws.Range("B68").Value = "คิดภาษีอย่างไร"

$ws.Range("A69").Value = "law-calculate"
$ws.Range("B69").Value = "คิดภาษียังไง"

$ws.Range("A70").Value = "law-calculate"
$ws.Range("B70").Value = "คำนวนภาษีอย่างไร"

$ws.Range("A71").Value = "law-calculate"
$ws.Range("B71").Value = "คำนวนภาษียังไง"

$ws.Range("A72").Value = "law-calculate"
$ws.Range("B72").Value = "คิดภาษีได้ไง"

$ws.Range("A73").Value = "law-calculate"
$ws.Range("B73").Value = "คิดภาษีจากไหน"

$ws.Range("A74").Value = "law-time"
$ws.Range("B74").Value = "กฎหมายกำหนดให้บุคคลต้องทำการยื่นเสียภาษีในช่วงเดือนใด"

$ws.Range("A75").Value = "law-time"
$ws.Range("B75").Value = "กำหนดการยื่นเสียภาษี"

$ws.Range("A76").Value = "law-time"
$ws.Range("B76").Value = "ช่วงเวลาในการเสียภาษี"

$ws.Range("A77").Value = "law-time"
$ws.Range("B77").Value = "เสียภาษีตอนไหน"

$ws.Range("A78").Value = "law-time"
$ws.Range("B78").Value = "ต้องเสียภาษีตอนไหน"

$ws.Range("A79").Value = "law-time"
$ws.Range("B79").Value = "ต้องเสียภาษีช่วงไหน"

$ws.Range("A80").Value = "law-time"
$ws.Range("B80").Value = "เวลาเสียภาษี"

$ws.Range("A81").Value = "law-time"
$ws.Range("B81").Value = "กำหนดการเสียภาษี"

$ws.Range("A82").Value = "law-time"
$ws.Range("B82").Value = "เสียภาษีเดือนไหน"

$ws.Range("A83").Value = "law-place"
$ws.Range("B83").Value = "ถ้าต้องการที่จะเสียภาษี สามารถยื่นเสียภาษีได้ที่ไหนบ้าง "

$ws.Range("A84").Value = "law-place"
$ws.Range("B84").Value = "เสียภาษีที่ไหน"

$ws.Range("A85").Value = "law-place"
$ws.Range("B85").Value = "ทำการเสียภาษีที่ไหนได้บ้าง"

$ws.Range("A86").Value = "law-place"
$ws.Range("B86").Value = "ต้องเสียภาษีที่ไหน"

$ws.Range("A87").Value = "law-place"
$ws.Range("B87").Value = "ต้องการยื่นเสียภาษี"

$ws.Range("A88").Value = "law-place"
$ws.Range("B88").Value = "สถานที่จ่ายภาษี"

$ws.Range("A89").Value = "law-place"
$ws.Range("B89").Value = "สถานที่เสียภาษี"

$ws.Range("A90").Value = "law-salary"
$ws.Range("B90").Value = "บุคคลธรรมดาต้องมีเงินเดือนเท่าไร ถึงต้องยื่นภาษี "

$ws.Range("A91").Value = "law-salary"
$ws.Range("B91").Value = "เงินดือนเท่าไหร่ต้องเสียภาษี"

$ws.Range("A92").Value = "law-salary"
$ws.Range("B92").Value = "เงินเดือนกี่บาทต้องเสียภาษี"

$ws.Range("A93").Value = "law-salary"
$ws.Range("B93").Value = "ทำงานเงินเดือนเท่าไหร่ต้องเสียภาษี"

$ws.Range("A94").Value = "law-salary"
$ws.Range("B94").Value = "ยื่นเสียภาษีตอนเงินเดือนกี่บาท"

$excel.CutCopyMode = $false

# Restore the workbook-level view state (scroll position / active cell)
# to roughly where the author left off after adding this data.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 47
$win.ScrollColumn = 1
$ws.Range("E64").Select()

